# "Added log for pvp module"
#
# - Property1-match (sheet1): rename two headers
#     MatchKilledHero    -> MatchOpponentK
#     MatchBeKilledHero  -> MatchOpponentD
#   and append three new trailing columns (X, Y, Z) mirroring column W
#   (MatchOpponentHeroStar3) with new headers MatchOpponentHeroHP1/2/3,
#   extending the TRUE/FALSE list validation to cover them too.
# - Property1-match becomes the active sheet/tab with selection on Z4.
# - record1-match is no longer the active tab; its selection moves to I15.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Property1-match")
$ws2 = $wb.Worksheets.Item("record1-match")

# --- Rename the two headers on sheet1 ---
$ws1.Range("E1").Value = "MatchOpponentK"
$ws1.Range("F1").Value = "MatchOpponentD"

# --- Add new columns X, Y, Z mirroring column W (header + data rows) ---
$null = $ws1.Range("W1:W9").Copy()
$null = $ws1.Range("X1:Z9").PasteSpecial(-4122)

$ws1.Range("X1").Value = "MatchOpponentHeroHP1"
$ws1.Range("Y1").Value = "MatchOpponentHeroHP2"
$ws1.Range("Z1").Value = "MatchOpponentHeroHP3"

for ($r = 2; $r -le 9; $r++) {
    $wv = $ws1.Cells.Item($r, 23).Value()
    $ws1.Cells.Item($r, 24).Value = $wv
    $ws1.Cells.Item($r, 25).Value = $wv
    $ws1.Cells.Item($r, 26).Value = $wv
}

# Column G (7) widens slightly (19.125 -> ~22.875) as part of the layout change
$ws1.Columns.Item(7).ColumnWidth = 22.14

# --- Extend the TRUE/FALSE list validation to the new columns ---
$null = $ws1.Range("B7:W9").Validation.Delete()
$null = $ws1.Range("B7:Z9").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# --- record1-match: no longer the active tab; selection moves to I15 ---
$null = $ws2.Range("I15").Select()

# --- Property1-match becomes the active sheet/tab; selection on Z4 ---
$ws1.Activate()
$null = $ws1.Range("Z4").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 21
$win.ScrollRow = 1
